# Daily attendance processing - 2025-12-03 09:55:06
# Normalizes the "Recorded By" column (G) so that "System" (the
# automated recorder) is listed first, followed by the remaining
# recorder names/emails in reverse of their original order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$lastCol = $ws.UsedRange.Columns.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $raw = $cell.Value2

    if ($raw -eq $null) { continue }

    $text = [string]$raw
    if ($text.IndexOf(",") -lt 0) { continue }

    $parts = $text.Split(",")
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    if ($trimmed -notcontains "System") { continue }

    $reversed = @()
    for ($i = $trimmed.Count - 1; $i -ge 0; $i--) {
        $reversed += $trimmed[$i]
    }

    $cell.Value2 = [string]::Join(", ", $reversed)
}
